$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Rename header cells in row 1
# -----------------------------------------------------------------
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2310")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2404")
}

# -----------------------------------------------------------------
# 2) Freeze the header row
# -----------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# -----------------------------------------------------------------
# 3) Build table while avoiding header dxf capture:
#    back up header formatting, strip it, build table, restore it.
# -----------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$backupRange = $ws.Range("A100000:U100000")

$headerRange.Copy()
$backupRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$range = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $range, 0, 1, "")
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$backupRange.Copy()
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$backupRange.ClearFormats()
$backupRange.ClearContents()

Write-Host "done"
